$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 9, pushing the existing rows 9-11 down to 10-12
$ws.Range("A9:R9").EntireRow.Insert()

# Populate the newly inserted row 9 with the new weekly record
$ws.Range("A9").Value2 = 7
$ws.Range("B9").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value2 = "Ñuble"
$ws.Range("D9").Value2 = 44880
$ws.Range("E9").Value2 = 16
$ws.Range("F9").Value2 = 100114007
$ws.Range("G9").Value2 = "Jengibre"
$ws.Range("H9").Value2 = "Sin especificar"
$ws.Range("I9").Value2 = "Primera"
$ws.Range("J9").Value2 = 30
$ws.Range("K9").Value2 = 17000
$ws.Range("L9").Value2 = 17000
$ws.Range("M9").Value2 = 17000
$ws.Range("N9").Value2 = "$/caja 13 kilos"
$ws.Range("O9").Value2 = "Perú"
$ws.Range("P9").Value2 = 1308
$ws.Range("Q9").Value2 = 13
$ws.Range("R9").Value2 = "Hortaliza"
